$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1536.3
$ws.Range("I6").Value = 1851.5714
$ws.Range("J6").Value = 800.6667
$ws.Range("K6").Value = 5554.7142
$ws.Range("L6").Value = 2402.0001
$ws.Range("M6").Value = -5442.7142
$ws.Range("N6").Value = -2626.0001

$ws.Range("H12").Value = 249.75
$ws.Range("I12").Value = 264
$ws.Range("K12").Value = 264
$ws.Range("M12").Value = -94

$ws.Range("H17").Value = 3394.7307
$ws.Range("J17").Value = 3417.52
$ws.Range("L17").Value = 10252.56
$ws.Range("N17").Value = -10588.56

$ws.Range("H28").Value = 101454.4
$ws.Range("I28").Value = 112532.664
$ws.Range("K28").Value = 112532.664
$ws.Range("M28").Value = -112047.664

$ws.Range("H49").Value = 2869
$ws.Range("J49").Value = 2869
$ws.Range("L49").Value = 8607
$ws.Range("N49").Value = -8879

$ws.Range("H86").Value = 321433630
$ws.Range("I86").Value = 333338780
$ws.Range("K86").Value = 333338780
$ws.Range("M86").Value = -333337657

$ws.Range("H89").Value = 321433630
$ws.Range("I89").Value = 333338780
$ws.Range("K89").Value = 1666693900
$ws.Range("M89").Value = -1666688284

$ws.Range("H98").Value = 3542.5715
$ws.Range("I98").Value = 1666.6666
$ws.Range("K98").Value = 1666.6666
$ws.Range("M98").Value = -168.6666

$ws.Range("H112").Value = 1284963.8
$ws.Range("J112").Value = 1336262.2
$ws.Range("L112").Value = 4008786.6
$ws.Range("N112").Value = -4011002.6

$ws.Range("H122").Value = 3542.5715
$ws.Range("I122").Value = 1666.6666
$ws.Range("K122").Value = 4999.9998
$ws.Range("M122").Value = -2549.9998

$ws.Range("H131").Value = 4680.727
$ws.Range("J131").Value = 5224.75
$ws.Range("L131").Value = 15674.25
$ws.Range("N131").Value = -25754.25

$ws.Range("H135").Value = 2054.5881
$ws.Range("I135").Value = 1554.5454
$ws.Range("J135").Value = 2971.3333
$ws.Range("K135").Value = 13990.9086
$ws.Range("L135").Value = 26741.9997
$ws.Range("M135").Value = -11455.9086
$ws.Range("N135").Value = -31811.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 198.93333
$ws.Range("I5").Value = 212.42857
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 212.42857
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = -100.42857
$ws.Range("N5").Value = -234

$ws.Range("H32").Value = 17490152
$ws.Range("I32").Value = 20725442
$ws.Range("K32").Value = 20725442
$ws.Range("M32").Value = -20725155

$ws.Range("H61").Value = 2540.2896
$ws.Range("I61").Value = 2308.2273
$ws.Range("K61").Value = 2308.2273
$ws.Range("M61").Value = -2096.2273

$ws.Range("H108").Value = 93179.39999999999
$ws.Range("J108").Value = 93179.39999999999
$ws.Range("L108").Value = 93179.39999999999
$ws.Range("N108").Value = -100859.4

$ws.Range("H132").Value = 5249.1333
$ws.Range("I132").Value = 5393.9165
$ws.Range("J132").Value = 4670
$ws.Range("K132").Value = 16181.7495
$ws.Range("L132").Value = 14010
$ws.Range("M132").Value = -13651.7495
$ws.Range("N132").Value = -19070

$ws.Range("H136").Value = 2540.2896
$ws.Range("I136").Value = 2308.2273
$ws.Range("K136").Value = 6924.6819
$ws.Range("M136").Value = -4374.6819

$ws.Range("H139").Value = 75282.14
$ws.Range("I139").Value = 69990
$ws.Range("K139").Value = 69990
$ws.Range("M139").Value = -64850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 198.93333
$ws.Range("I4").Value = 212.42857
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 212.42857
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = -97.42857000000001
$ws.Range("N4").Value = -240

$ws.Range("H11").Value = 230.5
$ws.Range("I11").Value = 240.66667
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 240.66667
$ws.Range("L11").Value = 200
$ws.Range("M11").Value = -100.66667
$ws.Range("N11").Value = -480

$ws.Range("H63").Value = 100264.5
$ws.Range("J63").Value = 100264.5
$ws.Range("L63").Value = 100264.5
$ws.Range("N63").Value = -101636.5

$ws.Range("H66").Value = 100264.5
$ws.Range("J66").Value = 100264.5
$ws.Range("L66").Value = 300793.5
$ws.Range("N66").Value = -307657.5

$ws.Range("H82").Value = 59471.445
$ws.Range("I82").Value = 35851
$ws.Range("J82").Value = 88997
$ws.Range("K82").Value = 35851
$ws.Range("L82").Value = 88997
$ws.Range("M82").Value = -35468
$ws.Range("N82").Value = -89763

$ws.Range("H85").Value = 59471.445
$ws.Range("I85").Value = 35851
$ws.Range("J85").Value = 88997
$ws.Range("K85").Value = 35851
$ws.Range("L85").Value = 88997
$ws.Range("M85").Value = -34525
$ws.Range("N85").Value = -91649

$ws.Range("H92").Value = 99999.5
$ws.Range("J92").Value = 99999
$ws.Range("L92").Value = 99999
$ws.Range("N92").Value = -104991

$ws.Range("H99").Value = 2687.8
$ws.Range("I99").Value = 2013.0667
$ws.Range("K99").Value = 2013.0667
$ws.Range("M99").Value = -515.0667000000001

$ws.Range("H134").Value = 3403809.2
$ws.Range("I134").Value = 3970166.2
$ws.Range("K134").Value = 11910498.6
$ws.Range("M134").Value = -11907963.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 8757143
$ws.Range("I6").Value = 4216666.5
$ws.Range("K6").Value = 4216666.5
$ws.Range("M6").Value = -4216553.5

$ws.Range("H7").Value = 35.434784
$ws.Range("I7").Value = 32.35
$ws.Range("K7").Value = 32.35
$ws.Range("M7").Value = 80.65000000000001

$ws.Range("H22").Value = 310.5625
$ws.Range("I22").Value = 295.3846
$ws.Range("J22").Value = 376.33334
$ws.Range("K22").Value = 295.3846
$ws.Range("L22").Value = 376.33334
$ws.Range("M22").Value = 54.61540000000002
$ws.Range("N22").Value = -1076.33334

$ws.Range("H31").Value = 5110.5405
$ws.Range("I31").Value = 2258.2727
$ws.Range("J31").Value = 6317.269
$ws.Range("K31").Value = 2258.2727
$ws.Range("L31").Value = 6317.269
$ws.Range("M31").Value = -1963.2727
$ws.Range("N31").Value = -6907.269

$ws.Range("H34").Value = 5110.5405
$ws.Range("I34").Value = 2258.2727
$ws.Range("J34").Value = 6317.269
$ws.Range("K34").Value = 2258.2727
$ws.Range("L34").Value = 6317.269
$ws.Range("M34").Value = -2056.2727
$ws.Range("N34").Value = -6721.269

$ws.Range("H58").Value = 2703.3508
$ws.Range("I58").Value = 2282.2708
$ws.Range("J58").Value = 4949.1113
$ws.Range("K58").Value = 2282.2708
$ws.Range("L58").Value = 4949.1113
$ws.Range("M58").Value = -2079.2708
$ws.Range("N58").Value = -5355.1113

$ws.Range("H132").Value = 3456.5
$ws.Range("I132").Value = 3079.7036
$ws.Range("K132").Value = 9239.110799999999
$ws.Range("M132").Value = -6709.110799999999

$ws.Range("H134").Value = 1905.6897
$ws.Range("I134").Value = 1909.4642
$ws.Range("K134").Value = 5728.392599999999
$ws.Range("M134").Value = -3193.392599999999

$ws.Range("H136").Value = 2703.3508
$ws.Range("I136").Value = 2282.2708
$ws.Range("J136").Value = 4949.1113
$ws.Range("K136").Value = 6846.812399999999
$ws.Range("L136").Value = 14847.3339
$ws.Range("M136").Value = -4296.812399999999
$ws.Range("N136").Value = -19947.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2715.9167
$ws.Range("I5").Value = 2024.75
$ws.Range("J5").Value = 3061.5
$ws.Range("K5").Value = 6074.25
$ws.Range("L5").Value = 9184.5
$ws.Range("M5").Value = -5962.25
$ws.Range("N5").Value = -9408.5

$ws.Range("H9").Value = 5003125

$ws.Range("H32").Value = 1996.6666

$ws.Range("H132").Value = 456916.9
$ws.Range("I132").Value = 1743.7273
$ws.Range("J132").Value = 912090.0600000001
$ws.Range("K132").Value = 15693.5457
$ws.Range("L132").Value = 8208810.540000001
$ws.Range("M132").Value = -13163.5457
$ws.Range("N132").Value = -8213870.540000001

$ws.Range("H135").Value = 2715.9167
$ws.Range("I135").Value = 2024.75
$ws.Range("J135").Value = 3061.5
$ws.Range("K135").Value = 18222.75
$ws.Range("L135").Value = 27553.5
$ws.Range("M135").Value = -15687.75
$ws.Range("N135").Value = -32623.5

$ws.Range("H140").Value = 3384.3845
$ws.Range("I140").Value = 1499.5
$ws.Range("K140").Value = 4498.5
$ws.Range("M140").Value = 681.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 8727.608
$ws.Range("I2").Value = 35.058823
$ws.Range("K2").Value = 35.058823
$ws.Range("M2").Value = 77.94117700000001

$ws.Range("H102").Value = 2359.0715
$ws.Range("I102").Value = 2359.0715
$ws.Range("K102").Value = 2359.0715
$ws.Range("M102").Value = -737.0715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8704.632
$ws.Range("I7").Value = 8800
$ws.Range("J7").Value = 8498
$ws.Range("K7").Value = 8800
$ws.Range("L7").Value = 8498
$ws.Range("M7").Value = -8688
$ws.Range("N7").Value = -8722

$ws.Range("H16").Value = 4893.25
$ws.Range("I16").Value = 4592.2856
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 4592.2856
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = -4422.2856
$ws.Range("N16").Value = -7340

$ws.Range("H46").Value = 8795.85
$ws.Range("I46").Value = 375
$ws.Range("K46").Value = 375
$ws.Range("M46").Value = -187

$ws.Range("H122").Value = 15328.72
$ws.Range("I122").Value = 15959.368
$ws.Range("K122").Value = 47878.104
$ws.Range("M122").Value = -45428.104

$ws.Range("H126").Value = 8704.632
$ws.Range("I126").Value = 8800
$ws.Range("J126").Value = 8498
$ws.Range("K126").Value = 26400
$ws.Range("L126").Value = 25494
$ws.Range("M126").Value = -23930
$ws.Range("N126").Value = -30434

$ws.Range("H132").Value = 6874.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 6874.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 20623.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -25683.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 24980
$ws.Range("J49").Value = 24980
$ws.Range("L49").Value = 24980
$ws.Range("N49").Value = -25440

$ws.Range("H113").Value = 797.5
$ws.Range("I113").Value = 593
$ws.Range("K113").Value = 1779
$ws.Range("M113").Value = 391

$ws.Range("H132").Value = 2332.2188
$ws.Range("I132").Value = 2332.3845
$ws.Range("K132").Value = 6997.1535
$ws.Range("M132").Value = -4467.1535

Write-Output "Applied all changes"